$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the cell that previously read "Durise_600ml_6timer" -> "Diurese_600ml_6timer"
$ws.Range("T1").Value = "Diurese_600ml_6timer"

# Update the current selection to the whole second row (A2:XFD2), active cell A2
$ws.Rows("2:2").Select()
